# Update column F (dSF) values to reflect repulled/pushed data and mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 6
    4  = -4
    6  = -9
    7  = 3
    8  = -4
    9  = -1
    10 = -10
    11 = 0
    12 = -1
    13 = 1
    14 = 0
    15 = -1
    19 = -5
    20 = -1
    22 = -2
    24 = -2
    25 = -3
    26 = 1
    27 = -2
    28 = 2
    29 = 3
    30 = -2
    33 = 1
    34 = -4
    35 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
